$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 4-25 (weekly rolling data window) ---
$ws.Range("D4").Value = 44602
$ws.Range("J4").Value = 12000

$ws.Range("D5").Value = 44602
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 6000
$ws.Range("K5").Value = 2500
$ws.Range("L5").Value = 2500
$ws.Range("M5").Value = 2500
$ws.Range("O5").Value = "Provincia de Chacabuco"
$ws.Range("P5").Value = 25

$ws.Range("D6").Value = 44161
$ws.Range("J6").Value = 7000

$ws.Range("D7").Value = 44600
$ws.Range("J7").Value = 1300
$ws.Range("K7").Value = 3500
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = 3808
$ws.Range("O7").Value = "Región Metropolitana"
$ws.Range("P7").Value = 38

$ws.Range("D8").Value = 44189
$ws.Range("J8").Value = 16000
$ws.Range("K8").Value = 3000
$ws.Range("M8").Value = 3000
$ws.Range("P8").Value = 30

$ws.Range("D9").Value = 44187

$ws.Range("D10").Value = 44209
$ws.Range("J10").Value = 7000
$ws.Range("K10").Value = 2500
$ws.Range("M10").Value = 2750
$ws.Range("O10").Value = "Provincia de Chacabuco"
$ws.Range("P10").Value = 28

$ws.Range("D11").Value = 44181
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 12000
$ws.Range("K11").Value = 3000
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = 3000
$ws.Range("O11").Value = "Provincia de Chacabuco"
$ws.Range("P11").Value = 30

$ws.Range("D12").Value = 44245
$ws.Range("J12").Value = 9000
$ws.Range("O12").Value = "Región Metropolitana"

$ws.Range("D13").Value = 44245
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 5000
$ws.Range("K13").Value = 2500
$ws.Range("L13").Value = 2500
$ws.Range("M13").Value = 2500
$ws.Range("O13").Value = "Región Metropolitana"
$ws.Range("P13").Value = 25

$ws.Range("D14").Value = 44230
$ws.Range("J14").Value = 16000

$ws.Range("D15").Value = 44159
$ws.Range("J15").Value = 7000

$ws.Range("D16").Value = 44204
$ws.Range("J16").Value = 7000
$ws.Range("K16").Value = 3000
$ws.Range("M16").Value = 3000
$ws.Range("P16").Value = 30

$ws.Range("D17").Value = 44186
$ws.Range("J17").Value = 10000

$ws.Range("D18").Value = 44210
$ws.Range("J18").Value = 8800
$ws.Range("K18").Value = 2500
$ws.Range("M18").Value = 2750
$ws.Range("P18").Value = 28

$ws.Range("D19").Value = 44188
$ws.Range("J19").Value = 12000

$ws.Range("D20").Value = 44166
$ws.Range("J20").Value = 7000

$ws.Range("D21").Value = 44162
$ws.Range("J21").Value = 7000

$ws.Range("D22").Value = 44231
$ws.Range("J22").Value = 12000

$ws.Range("D23").Value = 44232
$ws.Range("J23").Value = 16000

$ws.Range("D24").Value = 44229
$ws.Range("J24").Value = 16000

$ws.Range("D25").Value = 44214

# --- Append two brand-new rows (26-27) ---
$ws.Range("A26").Value = 6
$ws.Range("B26").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C26").Value = "Metropolitana"
$ws.Range("D26").Value = 44167
$ws.Range("D26").NumberFormat = $ws.Range("D25").NumberFormat
$ws.Range("E26").Value = 13
$ws.Range("F26").Value = 300000001
$ws.Range("G26").Value = "Rabanito"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 7000
$ws.Range("K26").Value = 3000
$ws.Range("L26").Value = 3000
$ws.Range("M26").Value = 3000
$ws.Range("N26").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O26").Value = "Provincia de Chacabuco"
$ws.Range("P26").Value = 30
$ws.Range("Q26").Value = 100
$ws.Range("R26").Value = "Hortaliza"

$ws.Range("A27").Value = 6
$ws.Range("B27").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C27").Value = "Metropolitana"
$ws.Range("D27").Value = 44160
$ws.Range("D27").NumberFormat = $ws.Range("D25").NumberFormat
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = 300000001
$ws.Range("G27").Value = "Rabanito"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 7000
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = 3000
$ws.Range("N27").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O27").Value = "Provincia de Chacabuco"
$ws.Range("P27").Value = 30
$ws.Range("Q27").Value = 100
$ws.Range("R27").Value = "Hortaliza"

Write-Host "Edit applied."
